# Add "Save" column (H) to the s_vals sheet, as in the commit
# "add save column in s_vals sheets"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing header cell (G1, style index 1:
# bold font, thin border, centered/top aligned) onto the new header
# cell H1 so the new column matches the look of the other headers
# without creating a duplicate style entry.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Header text for the new column
$ws.Range("H1").Value = "Save"

# Data values for the new column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
